$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 631.74194
$ws.Range("J19").Value = 717.0526
$ws.Range("L19").Value = 717.0526
$ws.Range("N19").Value = -1067.0526
$ws.Range("H33").Value = 235.25
$ws.Range("I33").Value = 238.45454
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 238.45454
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -9.454540000000009
$ws.Range("N33").Value = -658
$ws.Range("H111").Value = 1842.3636
$ws.Range("I111").Value = 1586.8
$ws.Range("J111").Value = 2055.3333
$ws.Range("K111").Value = 4760.4
$ws.Range("L111").Value = 6165.999899999999
$ws.Range("M111").Value = -1693.4
$ws.Range("N111").Value = -12299.9999
$ws.Range("H116").Value = 2743.75
$ws.Range("I116").Value = 2942
$ws.Range("J116").Value = 1356
$ws.Range("K116").Value = 2942
$ws.Range("L116").Value = 1356
$ws.Range("M116").Value = 500
$ws.Range("N116").Value = -8240
$ws.Range("H129").Value = 1295
$ws.Range("I129").Value = 442
$ws.Range("J129").Value = 1636.2
$ws.Range("K129").Value = 1326
$ws.Range("L129").Value = 4908.6
$ws.Range("M129").Value = 3674
$ws.Range("N129").Value = -14908.6
$ws.Range("H131").Value = 11541.85
$ws.Range("I131").Value = 3444.75
$ws.Range("J131").Value = 23687.5
$ws.Range("K131").Value = 10334.25
$ws.Range("L131").Value = 71062.5
$ws.Range("M131").Value = -5294.25
$ws.Range("N131").Value = -81142.5
$ws.Range("H138").Value = 6369921.5
$ws.Range("I138").Value = 2152951.5
$ws.Range("J138").Value = 8623819
$ws.Range("K138").Value = 6458854.5
$ws.Range("L138").Value = 25871457
$ws.Range("M138").Value = -6453714.5
$ws.Range("N138").Value = -25881737

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 963.6
$ws.Range("I45").Value = 948
$ws.Range("K45").Value = 948
$ws.Range("M45").Value = -571
$ws.Range("H97").Value = 14992.714
$ws.Range("I97").Value = 17358.166
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 17358.166
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -16862.166
$ws.Range("N97").Value = -1792
$ws.Range("H132").Value = 2930.0637
$ws.Range("I132").Value = 2319
$ws.Range("J132").Value = 5191
$ws.Range("K132").Value = 6957
$ws.Range("L132").Value = 15573
$ws.Range("M132").Value = -4427
$ws.Range("N132").Value = -20633

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3075.487
$ws.Range("I134").Value = 2110.0667
$ws.Range("J134").Value = 6293.5557
$ws.Range("K134").Value = 6330.2001
$ws.Range("L134").Value = 18880.6671
$ws.Range("M134").Value = -3795.2001
$ws.Range("N134").Value = -23950.6671

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1310.7391
$ws.Range("I16").Value = 1062.5834
$ws.Range("J16").Value = 1581.4546
$ws.Range("K16").Value = 1062.5834
$ws.Range("L16").Value = 1581.4546
$ws.Range("M16").Value = -775.5834
$ws.Range("N16").Value = -2155.4546
$ws.Range("H17").Value = 49504.5
$ws.Range("I17").Value = 50000
$ws.Range("J17").Value = 49009
$ws.Range("K17").Value = 50000
$ws.Range("L17").Value = 49009
$ws.Range("M17").Value = -49826
$ws.Range("N17").Value = -49357
$ws.Range("H86").Value = 26317512
$ws.Range("I86").Value = 38462716
$ws.Range("K86").Value = 38462716
$ws.Range("M86").Value = -38461593
$ws.Range("H89").Value = 26317512
$ws.Range("I89").Value = 38462716
$ws.Range("K89").Value = 192313580
$ws.Range("M89").Value = -192307964
$ws.Range("H113").Value = 1310.7391
$ws.Range("I113").Value = 1062.5834
$ws.Range("J113").Value = 1581.4546
$ws.Range("K113").Value = 1062.5834
$ws.Range("L113").Value = 1581.4546
$ws.Range("M113").Value = 1107.4166
$ws.Range("N113").Value = -5921.4546
$ws.Range("H122").Value = 2511.8635
$ws.Range("I122").Value = 1518.75
$ws.Range("J122").Value = 3079.3572
$ws.Range("K122").Value = 4556.25
$ws.Range("L122").Value = 9238.071599999999
$ws.Range("M122").Value = -2106.25
$ws.Range("N122").Value = -14138.0716
$ws.Range("H132").Value = 11907922
$ws.Range("I132").Value = 27779628
$ws.Range("J132").Value = 4143.875
$ws.Range("K132").Value = 83338884
$ws.Range("L132").Value = 12431.625
$ws.Range("M132").Value = -83336354
$ws.Range("N132").Value = -17491.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H132").Value = 3761.7917
$ws.Range("I132").Value = 3699.2727
$ws.Range("J132").Value = 4449.5
$ws.Range("K132").Value = 11097.8181
$ws.Range("L132").Value = 13348.5
$ws.Range("M132").Value = -8567.8181
$ws.Range("N132").Value = -18408.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 229.96
$ws.Range("I55").Value = 210.75
$ws.Range("J55").Value = 306.8
$ws.Range("K55").Value = 210.75
$ws.Range("L55").Value = 306.8
$ws.Range("M55").Value = -37.75
$ws.Range("N55").Value = -652.8
$ws.Range("H68").Value = 1895.25
$ws.Range("I68").Value = 1893.3334
$ws.Range("J68").Value = 1901
$ws.Range("K68").Value = 1893.3334
$ws.Range("L68").Value = 1901
$ws.Range("M68").Value = -1144.3334
$ws.Range("N68").Value = -3399
$ws.Range("H71").Value = 1895.25
$ws.Range("I71").Value = 1893.3334
$ws.Range("J71").Value = 1901
$ws.Range("K71").Value = 9466.666999999999
$ws.Range("L71").Value = 9505
$ws.Range("M71").Value = -5722.666999999999
$ws.Range("N71").Value = -16993
$ws.Range("H127").Value = 48000
$ws.Range("J127").Value = 48000
$ws.Range("L127").Value = 48000
$ws.Range("N127").Value = -57920

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 12999.5
$ws.Range("I17").Value = 12999.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 12999.5
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -12827.5
$ws.Range("N17").ClearContents()
$ws.Range("H122").Value = 2751.8572
$ws.Range("I122").Value = 2710.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8131.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5681.5
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2290.7612
$ws.Range("I132").Value = 2125.4897
$ws.Range("J132").Value = 2740.6667
$ws.Range("K132").Value = 6376.4691
$ws.Range("L132").Value = 8222.000100000001
$ws.Range("M132").Value = -3846.4691
$ws.Range("N132").Value = -13282.0001
$ws.Range("H136").Value = 2863.2983
$ws.Range("I136").Value = 1066.475
$ws.Range("J136").Value = 7091.1177
$ws.Range("K136").Value = 3199.425
$ws.Range("L136").Value = 21273.3531
$ws.Range("M136").Value = -649.4249999999997
$ws.Range("N136").Value = -26373.3531
